# Glowing eyes genes - 3023336989 : translation sheet refresh
# 1) Keep the existing "Main_240327" sheet, but rename it "Old_240327" and
#    move it to the back.
# 2) Insert a brand-new sheet "Main_250630" in front of it, rebuilt from the
#    up to date translation export (adds the White / Noctol genes, drops the
#    ThingDef mote rows and the English-source column that nobody filled in).

$wb = $excel.ActiveWorkbook

# Grab the pre-existing worksheet BEFORE inserting anything else - worksheet
# handles in this host resolve by live index, so do this first and keep a
# handle that we re-fetch by index *after* the insert below.
$oldWs = $wb.Worksheets.Item(1)

# Insert the new sheet; Add() puts it immediately before the (currently)
# active sheet, i.e. at position 1.
$newWs = $wb.Worksheets.Add()

# Re-resolve the original sheet now that indices have shifted, then rename
# both tabs to match the refreshed workbook.
$oldWs = $wb.Worksheets.Item(2)
$oldWs.Name = "Old_240327"
$newWs.Name = "Main_250630"

$fontName = "맑은 고딕"

function Set-Cell($ws, $row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $text
    $c.Font.Name = $fontName
}

# ---- Header row ----
Set-Cell $newWs 1 1 "Class+Node [(Identifier (Key)]"
Set-Cell $newWs 1 2 "Class [Not chosen]"
Set-Cell $newWs 1 3 "Node [Not chosen]"
Set-Cell $newWs 1 4 "Required Mods [Not chosen]"
Set-Cell $newWs 1 5 "English [Source string]"
Set-Cell $newWs 1 6 "Korean (한국어) [Translation]"

# ---- Data rows (GeneDef label/description pairs) ----
Set-Cell $newWs 2 1 "GeneDef+GEG_Eyes_GlowingRed.label"
Set-Cell $newWs 2 2 "GeneDef"
Set-Cell $newWs 2 3 "GEG_Eyes_GlowingRed.label"
Set-Cell $newWs 2 6 "빛나는 붉은 눈"

Set-Cell $newWs 3 1 "GeneDef+GEG_Eyes_GlowingRed.description"
Set-Cell $newWs 3 2 "GeneDef"
Set-Cell $newWs 3 3 "GEG_Eyes_GlowingRed.description"
Set-Cell $newWs 3 6 "이 유전자의 보유자는 홍채에 생물발광성 '핏빛 붉은색' 색소를 가지고 있습니다."

Set-Cell $newWs 4 1 "GeneDef+GEG_Eyes_GlowingArchotech.label"
Set-Cell $newWs 4 2 "GeneDef"
Set-Cell $newWs 4 3 "GEG_Eyes_GlowingArchotech.label"
Set-Cell $newWs 4 6 "빛나는 초월공학 눈"

Set-Cell $newWs 5 1 "GeneDef+GEG_Eyes_GlowingArchotech.description"
Set-Cell $newWs 5 2 "GeneDef"
Set-Cell $newWs 5 3 "GEG_Eyes_GlowingArchotech.description"
Set-Cell $newWs 5 6 "이 유전자의 보유자는 홍채에 생물발광성 '초월공학 라임색' 색소를 가지고 있습니다."

Set-Cell $newWs 6 1 "GeneDef+GEG_Eyes_GlowingPurple.label"
Set-Cell $newWs 6 2 "GeneDef"
Set-Cell $newWs 6 3 "GEG_Eyes_GlowingPurple.label"
Set-Cell $newWs 6 6 "빛나는 보라색 눈"

Set-Cell $newWs 7 1 "GeneDef+GEG_Eyes_GlowingPurple.description"
Set-Cell $newWs 7 2 "GeneDef"
Set-Cell $newWs 7 3 "GEG_Eyes_GlowingPurple.description"
Set-Cell $newWs 7 6 "이 유전자의 보유자는 홍채에 생물발광성 보라색 색소를 가지고 있습니다."

Set-Cell $newWs 8 1 "GeneDef+GEG_Eyes_GlowingOrange.label"
Set-Cell $newWs 8 2 "GeneDef"
Set-Cell $newWs 8 3 "GEG_Eyes_GlowingOrange.label"
Set-Cell $newWs 8 6 "빛나는 주황색 눈"

Set-Cell $newWs 9 1 "GeneDef+GEG_Eyes_GlowingOrange.description"
Set-Cell $newWs 9 2 "GeneDef"
Set-Cell $newWs 9 3 "GEG_Eyes_GlowingOrange.description"
Set-Cell $newWs 9 6 "이 유전자의 보유자는 홍채에 생물발광성 '불타는 주황색' 색소를 가지고 있습니다."

Set-Cell $newWs 10 1 "GeneDef+GEG_Eyes_GlowingCyan.label"
Set-Cell $newWs 10 2 "GeneDef"
Set-Cell $newWs 10 3 "GEG_Eyes_GlowingCyan.label"
Set-Cell $newWs 10 6 "빛나는 청록색 눈"

Set-Cell $newWs 11 1 "GeneDef+GEG_Eyes_GlowingCyan.description"
Set-Cell $newWs 11 2 "GeneDef"
Set-Cell $newWs 11 3 "GEG_Eyes_GlowingCyan.description"
Set-Cell $newWs 11 6 "이 유전자의 보유자는 홍채에 생물발광성 '영혼빛 청록색' 색소를 가지고 있습니다."

Set-Cell $newWs 12 1 "GeneDef+GEG_Eyes_GlowingWhite.label"
Set-Cell $newWs 12 2 "GeneDef"
Set-Cell $newWs 12 3 "GEG_Eyes_GlowingWhite.label"
Set-Cell $newWs 12 6 "빛나는 백색 눈"

Set-Cell $newWs 13 1 "GeneDef+GEG_Eyes_GlowingWhite.description"
Set-Cell $newWs 13 2 "GeneDef"
Set-Cell $newWs 13 3 "GEG_Eyes_GlowingWhite.description"
Set-Cell $newWs 13 6 "이 유전자의 보유자는 홍채에 생물발광성 순백색 색소를 가지고 있습니다."

Set-Cell $newWs 14 1 "GeneDef+GEG_Eyes_GlowingNoctol.label"
Set-Cell $newWs 14 2 "GeneDef"
Set-Cell $newWs 14 3 "GEG_Eyes_GlowingNoctol.label"
Set-Cell $newWs 14 6 "빛나는 암귀 눈"

Set-Cell $newWs 15 1 "GeneDef+GEG_Eyes_GlowingNoctol.description"
Set-Cell $newWs 15 2 "GeneDef"
Set-Cell $newWs 15 3 "GEG_Eyes_GlowingNoctol.description"
Set-Cell $newWs 15 6 "이 유전자의 보유자는 홍채에 생물발광성 '암귀의 노란색' 색소를 가지고 있습니다."

# ---- Column widths (best-fit, approximated) ----
$newWs.Columns.Item(1).ColumnWidth = 56.28515625
$newWs.Columns.Item(2).ColumnWidth = 21
$newWs.Columns.Item(3).ColumnWidth = 44.85546875
$newWs.Columns.Item(4).ColumnWidth = 32.28515625
$newWs.Columns.Item(5).ColumnWidth = 25.28515625
$newWs.Columns.Item(6).ColumnWidth = 91.28515625

# ---- Selections: old sheet keeps its own cursor, new sheet becomes active ----
$oldWs.Activate()
$oldWs.Range("A47").Select()

$newWs.Activate()
$newWs.Range("E25").Select()
